$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells remain text (matching original inlineStr type)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "265.77"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "9"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.61"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "9"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.120"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "9"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06099"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "9"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.567"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "9"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.482"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "9"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.343"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "9"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8201"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "9"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01335"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "9"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1586"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "9"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08072"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "9"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03454"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "9"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03209"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "9"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09224"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "9"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.749"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "9"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001624"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "9"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04646"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "9"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006348"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "9"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.006141"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "9"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001069"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "9"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "9"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.726"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "9"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.263"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "9"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3314"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "9"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1243"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "9"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "9"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "9"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "9"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "9"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "9"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "9"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "9"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "9"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "9"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "9"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "9"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "9"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "9"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04601"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "9"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007003"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "9"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1116"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "9"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003452"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "9"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01170"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "9"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005953"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "9"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0009907"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "9"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "9"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.8030"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "9"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001124"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "9"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "9"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "9"
